$wb = $excel.ActiveWorkbook

# --- Step 1: rename the existing "总计" sheet out of the way (it becomes
# "2022-Q1" further below) so a freshly-added sheet can take the "总计"
# name, then add the new "总计" (totals) sheet at the end. ---
$oldTotals = $wb.Worksheets.Item("总计")
$oldTotals.Name = "2022-Q1__tmp"
$newTotals = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$newTotals.Name = "总计"

# Copy header style (bold + border) from the old totals sheet header row.
$oldTotals.Range("B1:D1").Copy()
$newTotals.Range("B1:D1").PasteSpecial(-4122)

# Copy the row-index column style (column A) too.
$oldTotals.Range("A2:A6").Copy()
$newTotals.Range("A2:A7").PasteSpecial(-4122)

# Header row for the new totals sheet.
$newTotals.Range("B1").Value = "日期"
$newTotals.Range("C1").Value = "持有数量(只)"
$newTotals.Range("D1").Value = "持有市值(亿元)"

# Data rows: new 2022-Q1 row on top, followed by the previous totals rows
# shifted down by one.
$totalsData = @(
    @(0, "2022-Q1", 5, 1.01),
    @(1, "2021-Q4", 1, 1),
    @(2, "2021-Q3", 1, 0.76),
    @(3, "2021-Q2", 1, 0.66),
    @(4, "2021-Q1", 1, 0.5600000000000001),
    @(5, "2020-Q4", 1, 0.61)
)

$r = 2
foreach ($row in $totalsData) {
    $newTotals.Cells.Item($r, 1).Value = $row[0]
    $newTotals.Cells.Item($r, 2).Value = $row[1]
    $newTotals.Cells.Item($r, 3).Value = $row[2]
    $newTotals.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# --- Step 2: turn the old "总计" sheet into the new "2022-Q1" fund-holdings
# sheet (same sheetId, renamed + new content/layout). ---
$q1 = $oldTotals
$q1.Name = "2022-Q1"


# Pull header + index-column formatting from an existing per-fund sheet
# (2021-Q4) which already has the exact target layout (columns B..H).
$template = $wb.Worksheets.Item("2021-Q4")
$template.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$q1.Range("A2:A6").PasteSpecial(-4122)

# Clear the old C2:D6 leftover content that falls outside the new header
# band, then rewrite header row.
$q1.Range("A1:H6").ClearContents()

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Column B (fund code, has leading zeros) and columns D, E, F, G (numeric-
# looking figures) are all stored as TEXT in the source data, so force text
# format before assigning them.
$q1.Range("B2:B6").NumberFormat = "@"
$q1.Range("D2:G6").NumberFormat = "@"

$fundData = @(
    @(0, "009693", "富国积极成长一年定期开放混合", "17.82", "98.74", "3.52", "0.6273", 4),
    @(1, "014269", "嘉实北交所精选两年定期混合A", "5.00", "87.43", "3.70", "0.1850", 8),
    @(2, "014273", "广发北交所精选两年定开混合A", "4.55", "52.69", "3.13", "0.1424", 10),
    @(3, "014274", "广发北交所精选两年定开混合C", "0.92", "52.69", "3.13", "0.0288", 10),
    @(4, "014270", "嘉实北交所精选两年定期混合C", "0.64", "87.43", "3.70", "0.0237", 8)
)

$r = 2
foreach ($row in $fundData) {
    $q1.Cells.Item($r, 1).Value = $row[0]
    $q1.Cells.Item($r, 2).Value = $row[1]
    $q1.Cells.Item($r, 3).Value = $row[2]
    $q1.Cells.Item($r, 4).Value = $row[3]
    $q1.Cells.Item($r, 5).Value = $row[4]
    $q1.Cells.Item($r, 6).Value = $row[5]
    $q1.Cells.Item($r, 7).Value = $row[6]
    $q1.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Restore the originally active sheet/tab (our edits above left "总计"
# selected as a side effect of being the last-touched sheet).
$wb.Worksheets.Item(1).Activate()
